$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "kekkk"
$ws.Range("D3").Value = "www"
$ws.Range("H3").Value = 22.0
$ws.Range("G8").Value = "location 22"
$ws.Range("H8").Value = "name 22"
